$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row 2 (FAPs/Rln3/Rxfp2/ECs), shifting row 3 (FAPs/Rln3/Rxfp2/MuSCs) up to row 2
$ws.Rows.Item(2).Delete()

# Update the recomputed TPM-derived values for the (now) row 2
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.026469
$ws.Range("N2").Value = 0.052938
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.030494308467
$ws.Range("R2").Value = 0.182965850802
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
